$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell changes from the match scorecard correction
# (batting/bowling data for both teams + totals were mixed up and are now fixed)

$ws.Range("A2").Value = 'Mohammad Rizwan'
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 'Caught'
$ws.Range("E2").Value = ' Mark Wood'
$ws.Range("J2").Value = 'Jason Roy'
$ws.Range("K2").Value = 67
$ws.Range("L2").Value = 22
$ws.Range("M2").Value = 'NOT OUT'
$ws.Range("N2").Value = ' '
$ws.Range("A3").Value = 'Babar Azam(C)'
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Adil Rashid'
$ws.Range("J3").Value = 'Jos Buttler'
$ws.Range("K3").Value = 27
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 'Bowled'
$ws.Range("N3").Value = ' Shaheen Afridi'
$ws.Range("A4").Value = 'Fakhar Zaman'
$ws.Range("B4").Value = 85
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 'LBW'
$ws.Range("E4").Value = ' Chris Woakes'
$ws.Range("J4").Value = 'Dawid Malan'
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 'Caught'
$ws.Range("N4").Value = ' Hasan Ali'
$ws.Range("A5").Value = 'Mohammad Hafeez'
$ws.Range("E5").Value = ' Adil Rashid'
$ws.Range("J5").Value = 'Jonny Bairstow'
$ws.Range("K5").Value = 17
$ws.Range("L5").Value = 5
$ws.Range("N5").Value = ' Hasan Ali'
$ws.Range("A6").Value = 'Shoaib Malik'
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 'Caught'
$ws.Range("E6").Value = ' Chris Jordan'
$ws.Range("J6").Value = 'Eoin Morgan(C)'
$ws.Range("L6").Value = 8
$ws.Range("M6").Value = 'Caught'
$ws.Range("N6").Value = ' Shadab Khan'
$ws.Range("A7").Value = 'Asif Ali'
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 'Caught'
$ws.Range("E7").Value = ' Mark Wood'
$ws.Range("J7").Value = 'Moeen Ali'
$ws.Range("K7").Value = 9
$ws.Range("M7").Value = 'LBW'
$ws.Range("N7").Value = ' Haris Rauf'
$ws.Range("A8").Value = 'Shadab Khan'
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 'LBW'
$ws.Range("E8").Value = ' Adil Rashid'
$ws.Range("J8").Value = 'Liam Livingstone'
$ws.Range("K8").Value = 63
$ws.Range("L8").Value = 20
$ws.Range("M8").Value = 'Caught'
$ws.Range("N8").Value = ' Shaheen Afridi'
$ws.Range("A9").Value = 'Imad Wasim'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 'LBW'
$ws.Range("E9").Value = ' Adil Rashid'
$ws.Range("J9").Value = 'Chris Woakes'
$ws.Range("K9").Value = 9
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = '* NOT OUT'
$ws.Range("N9").Value = ' '
$ws.Range("A10").Value = 'Hasan Ali'
$ws.Range("B10").Value = 59
$ws.Range("C10").Value = 19
$ws.Range("D10").Value = 'LBW'
$ws.Range("E10").Value = ' Adil Rashid'
$ws.Range("J10").Value = 'Chris Jordan'
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ' '
$ws.Range("A11").Value = 'Shaheen Afridi'
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 'LBW'
$ws.Range("E11").Value = ' Chris Jordan'
$ws.Range("J11").Value = 'Adil Rashid'
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ' '
$ws.Range("A12").Value = 'Haris Rauf'
$ws.Range("B12").Value = 13
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 'NOT OUT'
$ws.Range("E12").Value = ' '
$ws.Range("J12").Value = 'Mark Wood'
$ws.Range("A16").Value = 213
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '13.5'
$ws.Range("D16").Value = 83
$ws.Range("J16").Value = 215
$ws.Range("K16").Value = 6
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = '11.5'
$ws.Range("M16").Value = 71
$ws.Range("A21").Value = 'Liam Livingstone'
$ws.Range("C21").Value = 30
$ws.Range("E21").Value = 15
$ws.Range("J21").Value = 'Hasan Ali'
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = '2.0'
$ws.Range("L21").Value = 41
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 20.5
$ws.Range("A22").Value = 'Mark Wood'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '3.0'
$ws.Range("C22").Value = 53
$ws.Range("E22").Value = 17.67
$ws.Range("J22").Value = 'Imad Wasim'
$ws.Range("L22").Value = 33
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 16.5
$ws.Range("A23").Value = 'Adil Rashid'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '3.0'
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 13.33
$ws.Range("J23").Value = 'Shadab Khan'
$ws.Range("L23").Value = 33
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 16.5
$ws.Range("A24").Value = 'Chris Jordan'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '3.0'
$ws.Range("C24").Value = 44
$ws.Range("E24").Value = 14.67
$ws.Range("J24").Value = 'Haris Rauf'
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = '3.0'
$ws.Range("L24").Value = 66
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 22
$ws.Range("A25").Value = 'Chris Woakes'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '2.5'
$ws.Range("C25").Value = 46
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 18.4
$ws.Range("J25").Value = 'Shaheen Afridi'
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = '2.5'
$ws.Range("L25").Value = 42
$ws.Range("M25").Value = 2
$ws.Range("N25").Value = 16.8
